$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.314
$ws.Range("A12").Value = -21.626
$ws.Range("E14").Value = 16.876
$ws.Range("E26").Value = 16.414
$ws.Range("A27").Value = -21.86
$ws.Range("E31").Value = 16.258
$ws.Range("A32").Value = -21.644
$ws.Range("E35").Value = 16.63
$ws.Range("A36").Value = -20.339
$ws.Range("E37").Value = 16.737
$ws.Range("A38").Value = -19.741
$ws.Range("E45").Value = 16.851
$ws.Range("A46").Value = -21.864
$ws.Range("E52").Value = 16.817
$ws.Range("A54").Value = -21.764
$ws.Range("A55").Value = -22.247
$ws.Range("A56").Value = -22.01
$ws.Range("E57").Value = 16.643
$ws.Range("A67").Value = -21.565
$ws.Range("A69").Value = -21.721
$ws.Range("A72").Value = -21.445
$ws.Range("E81").Value = 16.638
$ws.Range("A83").Value = -21.877
$ws.Range("E83").Value = 16.554
$ws.Range("A86").Value = -22.093
$ws.Range("A91").Value = -21.652
$ws.Range("A93").Value = -21.259
$ws.Range("A99").Value = -20.43
$ws.Range("E100").Value = 16.568
$ws.Range("E102").Value = 16.669
